# Apply cryptos list update (prices & volume-1h changes, plus a few row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.854.11"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").Value = "2.579.37"
$ws.Range("E3").Value = "  -5.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.88"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.73"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.42%  "
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "3.027.94"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.44"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").Value = "61.739.62"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").Value = "2.569.56"
$ws.Range("E17").Value = "  -5.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.58"
$ws.Range("E18").Value = "  -4.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.55"
$ws.Range("E19").Value = "  -2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.95"
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.05"
$ws.Range("E21").Value = "  -5.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.495"
$ws.Range("E23").Value = "  -3.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.08"
$ws.Range("E24").Value = "  -2.05%  "
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.06"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.38"
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0839"
$ws.Range("E29").Value = "  -4.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.92"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("E31").Value = "  -5.00%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.37"
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.73"
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.21"
$ws.Range("E35").Value = "  -3.32%  "
$ws.Range("E36").Value = "  -4.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.81"
$ws.Range("E37").Value = "  +1.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "336.10"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("B39").Value = "SuiNetwork"
$ws.Range("C39").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.934"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.08"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.95"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.51"
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.52"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.142.94"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.605"
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.64"
$ws.Range("E48").Value = "  -4.94%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0547"
$ws.Range("E49").Value = "  -4.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0967"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0240"
$ws.Range("E51").Value = "  -1.82%  "
